$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shortage")

$ws.Range("B2").Value = 0.5934328645957652
$ws.Range("C2").Value = 0.00001364928244292737
$ws.Range("D2").Value = 0.0001104115643712736

$ws.Range("B3").Value = 15.37602248253422
$ws.Range("C3").Value = 0.1219970040973848
$ws.Range("D3").Value = 0.00002330505244000491

$ws.Range("B4").Value = 3.451610643570803
$ws.Range("C4").Value = 0.0001177553333192535
$ws.Range("D4").Value = 0.000001757868254029659

$ws.Range("B6").Value = 37.81914754926367

$ws.Range("B7").Value = 8.489644988493703
